# Automatische test-sync: 2025-06-26 23:58:50
# Append a new "Testmail #18" row to the Logs sheet and bump the
# corresponding Dashboard aggregate count.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# New row goes right after the current last data row (row 49 -> row 50).
$newRow = 50

$logs.Cells.Item($newRow, 1).Value = "Kun je dit product voor mij bestellen?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #18: Kun je dit product voor mij bestellen?"
$logs.Cells.Item($newRow, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($newRow, 5).Value = "Geachte klant,
Dank u wel voor uw interesse in ons product. Helaas kunnen wij op basis van deze e-mail geen bestelling voor u plaatsen. U kunt echter onze website bezoeken en het product zelf bestellen. Mocht u verdere vragen hebben of hulp nodig hebben bij het plaatsen van een bestelling, dan helpen wij u graag verder.
Met vriendelijke groet,
[Naam Bedrijf] E-mailassistent"
$logs.Cells.Item($newRow, 6).Value = "2025-06-26 23:58:23"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"

# Update the conditional-formatting ranges so they keep covering the
# newly added row (D/G/H/I, rows 2:49 -> 2:50).
$dFc = $logs.Range("D2:D49").FormatConditions
for ($i = 1; $i -le $dFc.Count; $i++) {
    $dFc.Item($i).ModifyAppliesToRange($logs.Range("D2:D50"))
}

$gFc = $logs.Range("G2:G49").FormatConditions
for ($i = 1; $i -le $gFc.Count; $i++) {
    $gFc.Item($i).ModifyAppliesToRange($logs.Range("G2:G50"))
}

$hFc = $logs.Range("H2:H49").FormatConditions
for ($i = 1; $i -le $hFc.Count; $i++) {
    $hFc.Item($i).ModifyAppliesToRange($logs.Range("H2:H50"))
}

$iFc = $logs.Range("I2:I49").FormatConditions
for ($i = 1; $i -le $iFc.Count; $i++) {
    $iFc.Item($i).ModifyAppliesToRange($logs.Range("I2:I50"))
}

# Bump the "Bestelling / Levering" tally on the Dashboard sheet (21 -> 22).
$dash.Cells.Item(2, 2).Value = 22
